$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REPORT")
$ws.Activate()

# Populate column H (testps*) first so their shared-string entries are
# created before the B/C/E strings, matching the target string order.
$ws.Range("H28").Value = "testps9023"
$ws.Range("H29").Value = "testps01"
$ws.Range("H30").Value = "magic_qq_appl"
$ws.Range("H31").Value = "testps0324"

$ws.Range("B28").Value = "Selection 1 (Date)"
$ws.Range("C28").Value = " "

# "mm-dd-yy" resolves to Excel's built-in date number format (id 14).
$ws.Range("D28").NumberFormat = "mm-dd-yy"
$ws.Range("D28").Value = (Get-Date -Year 2020 -Month 1 -Day 9 -Hour 0 -Minute 0 -Second 0)

$ws.Range("E28").Value = "Total Changes"
$ws.Range("G28").Value = 4

# D29 carries the same style as D28 (built-in date number format) but stays empty.
$ws.Range("D28").Copy()
$ws.Range("D29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$excel.Goto($ws.Range("A10"), $true)
$ws.Range("B29:G29").Select()
